$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.236.79"
$ws.Range("E2").Value = "  +0.07%  "

$ws.Range("D3").Value = "1.905.32"
$ws.Range("E3").Value = "  +0.46%  "

$ws.Range("D5").Value = "'306.48"
$ws.Range("E5").Value = "  -0.31%  "

$ws.Range("D7").Value = "'0.5360"
$ws.Range("E7").Value = "  +3.05%  "

$ws.Range("D8").Value = "'0.3813"
$ws.Range("E8").Value = "  +1.29%  "

$ws.Range("D9").Value = "'0.07299"
$ws.Range("E9").Value = "  +0.17%  "

$ws.Range("D10").Value = "'22.26"
$ws.Range("E10").Value = "  +5.07%  "

$ws.Range("D11").Value = "'0.9059"
$ws.Range("E11").Value = "  +0.60%  "

$ws.Range("D12").Value = "'0.08225"
$ws.Range("E12").Value = "  +0.60%  "

$ws.Range("D13").Value = "'95.77"
$ws.Range("E13").Value = "  -0.91%  "

$ws.Range("D14").Value = "'5.350"
$ws.Range("E14").Value = "  +1.27%  "

$ws.Range("D15").Value = "'1.002"
$ws.Range("E15").Value = "  +0.02%  "

$ws.Range("D16").Value = "'14.89"
$ws.Range("E16").Value = "  +2.39%  "

$ws.Range("D17").Value = "'0.000008670"
$ws.Range("E17").Value = "  +0.66%  "

$ws.Range("E18").Value = "  -0.02%  "

$ws.Range("D19").Value = "27.257.35"
$ws.Range("E19").Value = "  +0.03%  "

$ws.Range("D20").Value = "'5.052"
$ws.Range("E20").Value = "  -0.68%  "

$ws.Range("D21").Value = "1.057.59"
$ws.Range("E21").Value = "  -44.22%  "

$ws.Range("E22").Value = "  +0.74%  "

$ws.Range("D23").Value = "'6.524"
$ws.Range("E23").Value = "  +1.84%  "

$ws.Range("D24").Value = "'149.01"
$ws.Range("E24").Value = "  +0.99%  "

$ws.Range("D25").Value = "'2.292"
$ws.Range("E25").Value = "  -0.10%  "

$ws.Range("E26").Value = "  +0.95%  "

$ws.Range("D27").Value = "'1.745"
$ws.Range("E27").Value = "  -0.17%  "

$ws.Range("D28").Value = "'116.91"
$ws.Range("E28").Value = "  +1.58%  "

$ws.Range("D29").Value = "'4.821"
$ws.Range("E29").Value = "  -0.38%  "

$ws.Range("D30").Value = "'4.730"
$ws.Range("E30").Value = "  -4.55%  "

$ws.Range("D31").Value = "'0.09228"
$ws.Range("E31").Value = "  +0.08%  "

$ws.Range("D32").Value = "'0.8300"
$ws.Range("E32").Value = "  +4.39%  "

$ws.Range("D33").Value = "'0.05081"
$ws.Range("E33").Value = "  +0.96%  "

$ws.Range("D34").Value = "'1.219"
$ws.Range("E34").Value = "  -0.03%  "

$ws.Range("D35").Value = "'3.005"
$ws.Range("E35").Value = "  +2.09%  "

$ws.Range("D36").Value = "'3.340"
$ws.Range("E36").Value = "  -3.12%  "

$ws.Range("D37").Value = "'2.676"
$ws.Range("E37").Value = "  +3.16%  "

$ws.Range("D38").Value = "'0.5875"
$ws.Range("E38").Value = "  +3.64%  "

$ws.Range("D39").Value = "'0.02004"
$ws.Range("E39").Value = "  +0.89%  "

$ws.Range("D40").Value = "'1.077"
$ws.Range("E40").Value = "  +0.33%  "

$ws.Range("D41").Value = "'9.331"
$ws.Range("E41").Value = "  +4.19%  "

$ws.Range("D42").Value = "'6.641"
$ws.Range("E42").Value = "  +1.25%  "

$ws.Range("D43").Value = "'117.33"
$ws.Range("E43").Value = "  +1.57%  "

$ws.Range("D44").Value = "'0.5093"
$ws.Range("E44").Value = "  +4.04%  "

$ws.Range("D45").Value = "'0.1525"
$ws.Range("E45").Value = "  +0.57%  "

$ws.Range("D46").Value = "'1.001"
$ws.Range("E46").Value = "  -0.06%  "

$ws.Range("D47").Value = "'10.09"
$ws.Range("E47").Value = "  -0.14%  "

$ws.Range("D48").Value = "'1.643"
$ws.Range("E48").Value = "  +1.21%  "

$ws.Range("D49").Value = "'38.38"
$ws.Range("E49").Value = "  +0.52%  "

$ws.Range("E50").Value = "  +3.68%  "

$ws.Range("D51").Value = "'63.54"
$ws.Range("E51").Value = "  +0.13%  "
